# Fix Training Data Issue (#48)
# The "Date" column (BF) on this team-stats sheet held the sheet's own
# file-name-derived label ("5-16-2007-08") instead of the actual game
# date. Re-point every data row (BF2:BF31) at the correct ISO date
# "2008-05-16" (the stats were recorded the day after the file's nominal
# label, due to how NBA.com stats pages roll the date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctDate = "2008-05-16"

# Rows 2-31 all carry the stale date string in column BF.
$firstRow = 2
$lastRow  = 31
$rangeAddress = "BF" + $firstRow + ":BF" + $lastRow
$dateRange = $ws.Range($rangeAddress)

# Force the range to text formatting before writing the value so Excel
# stores the literal string "2008-05-16" instead of silently converting
# the ISO-looking text into a date serial number.
$dateRange.NumberFormat = "@"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Range("BF$row").Value = $correctDate
}

# Drop the temporary text-number-format again so the cells end up back
# on the sheet's normal (default) style, matching every other cell.
$dateRange.ClearFormats()
